# "Actualizacao das base de dados"
#
# The Chimanimani / Gorongosa partner-institution table (A:E, rows 2-10) is
# refreshed: the per-institution Masculino/Feminino/TOTAL counts themselves
# are unchanged, but the rows are re-sorted, so each institution (and its
# figures) now lands on a different row than before.
#
# Row 4 (SDAE SUSSUNDENGA) happens to stay put; every other row moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Region | Institution            | Masculino | Feminino | TOTAL
$rows = @(
    @(2, "CHIMANIMANI", "UNIZAMBEZE",           0,  1,  1),
    @(3, "CHIMANIMANI", "ITAM",                 1,  0,  1),
    @(4, "CHIMANIMANI", "SDAE SUSSUNDENGA",      1,  0,  1),
    @(5, "CHIMANIMANI", "ISPM",                 0,  1,  1),
    @(6, "CHIMANIMANI", "PARQUE DE CHIMANIMANI", 1,  0,  1),
    @(7, "CHIMANIMANI", "UCM",                  1,  0,  1),
    @(8, "CHIMANIMANI", "MICAIA",                3,  3,  6),
    @(9, "GORONGOSA",   "PARQUE DE GORONGOSA",  12, 11, 23),
    @(10,"GORONGOSA",   "GORONGOSA",             5,  4,  9)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]   # B: Instituicao
    $ws.Cells.Item($rowNum, 3).Value = $r[3]   # C: Masculino
    $ws.Cells.Item($rowNum, 4).Value = $r[4]   # D: Feminino
    $ws.Cells.Item($rowNum, 5).Value = $r[5]   # E: TOTAL
}
